# The source data pipeline was re-run with a fixed workflow: each sheet's
# table (Cutoff / Reaction_number) now starts at Cutoff-step 5 instead of 1,
# so the first 4 data rows are gone and every kept row's "Cutoff" (col B)
# value is shifted up by 4 while the "Reaction_number" (col C) values were
# recomputed. Column A (0-based index) stays 0..14. The table now ends at
# row 16 instead of row 20.

$wb = $excel.ActiveWorkbook

# New column B / C values (by sheet name) for data rows 2..16.
$newValues = @{
    "NBR" = @(
        @(5, 87),
        @(6, 85),
        @(7, 88),
        @(8, 88),
        @(9, 87),
        @(10, 86),
        @(11, 85),
        @(12, 85),
        @(13, 85),
        @(14, 85),
        @(15, 84),
        @(16, 84),
        @(17, 84),
        @(18, 84),
        @(19, 83)
    )
    "BAR" = @(
        @(5, 583),
        @(6, 581),
        @(7, 577),
        @(8, 577),
        @(9, 577),
        @(10, 572),
        @(11, 573),
        @(12, 572),
        @(13, 573),
        @(14, 573),
        @(15, 574),
        @(16, 572),
        @(17, 571),
        @(18, 572),
        @(19, 569)
    )
}

foreach ($ws in $wb.Worksheets) {
    $rows = $newValues[$ws.Name]
    if ($rows -eq $null) { continue }

    # Overwrite the 15 rows that survive (rows 2..16) with their new B/C values.
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $i + 2
        $pair = $rows[$i]
        $ws.Cells.Item($r, 2).Value = $pair[0]
        $ws.Cells.Item($r, 3).Value = $pair[1]
    }

    # Drop the now-stale trailing rows (old rows 17..20) so the used range
    # shrinks back down to A1:C16.
    $ws.Range("A17:A20").EntireRow.Delete()
}
